$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.78907395996592
$ws.Range("D2").Value = 9.460563024875897
$ws.Range("E2").Value = 14.5234661690076
$ws.Range("F2").Value = 34.12845656499626
$ws.Range("G2").Value = 36.0436097434874
$ws.Range("H2").Value = 16.27177381568231
$ws.Range("J2").Value = 10.63916787136682
$ws.Range("K2").Value = 18.25623686488967
$ws.Range("L2").Value = 9.837132795172533
$ws.Range("O2").Value = 25.70484035575479
$ws.Range("C3").Value = 13.73759570133807
$ws.Range("D3").Value = 9.408922606649135
$ws.Range("E3").Value = 14.51466324519637
$ws.Range("F3").Value = 34.29035512237162
$ws.Range("G3").Value = 36.29648887834801
$ws.Range("H3").Value = 16.35324828137597
$ws.Range("J3").Value = 10.66273249963868
$ws.Range("K3").Value = 17.60117637634765
$ws.Range("L3").Value = 9.847657745140191
$ws.Range("O3").Value = 25.85773459955056
$ws.Range("C4").Value = 13.70865341613171
$ws.Range("D4").Value = 9.378084436155707
$ws.Range("E4").Value = 14.51138196415795
$ws.Range("F4").Value = 34.39889614389239
$ws.Range("G4").Value = 36.46458138862744
$ws.Range("H4").Value = 16.40638446748228
$ws.Range("J4").Value = 10.67841396212094
$ws.Range("K4").Value = 17.18575067399386
$ws.Range("L4").Value = 9.85492716808427
$ws.Range("O4").Value = 25.95802771258034
$ws.Range("C5").Value = 13.69753846575042
$ws.Range("D5").Value = 9.365744698667514
$ws.Range("E5").Value = 14.51058096461735
$ws.Range("F5").Value = 34.44541726723241
$ws.Range("G5").Value = 36.53628909505326
$ws.Range("H5").Value = 16.42882032335166
$ws.Range("J5").Value = 10.68510956159761
$ws.Range("K5").Value = 17.01335021726958
$ws.Range("L5").Value = 9.858092916032504
$ws.Range("O5").Value = 26.00050883203947
$ws.Range("C6").Value = 13.6957341100128
$ws.Range("D6").Value = 9.363709651251014
$ws.Range("E6").Value = 14.51048039858657
$ws.Range("F6").Value = 34.45328014927596
$ws.Range("G6").Value = 36.54838943543956
$ws.Range("H6").Value = 16.43259304517129
$ws.Range("J6").Value = 10.68623980856178
$ws.Range("K6").Value = 16.98454199553545
$ws.Range("L6").Value = 9.858630884119306
$ws.Range("O6").Value = 26.00766001475527
$ws.Range("C7").Value = 13.70850075440194
$ws.Range("D7").Value = 9.377917087750426
$ws.Range("E7").Value = 14.51136898812405
$ws.Range("F7").Value = 34.39951428241525
$ws.Range("G7").Value = 36.46553549319795
$ws.Range("H7").Value = 16.40668387690287
$ws.Range("J7").Value = 10.67850302482075
$ws.Range("K7").Value = 17.18343792093052
$ws.Range("L7").Value = 9.854969038311436
$ws.Range("O7").Value = 25.95859410856103
$ws.Range("C8").Value = 13.77077663765619
$ws.Range("D8").Value = 9.442582168434427
$ws.Range("E8").Value = 14.51999147057066
$ws.Range("F8").Value = 34.18237898064548
$ws.Range("G8").Value = 36.12813126060702
$ws.Range("H8").Value = 16.29922077230788
$ws.Range("J8").Value = 10.64704146561079
$ws.Range("K8").Value = 18.03323528318746
$ws.Range("L8").Value = 9.84059456821168
$ws.Range("O8").Value = 25.75622533802693
$ws.Range("C9").Value = 13.91362076768727
$ws.Range("D9").Value = 9.575903330757098
$ws.Range("E9").Value = 14.55365016588498
$ws.Range("F9").Value = 33.82937573281237
$ws.Range("G9").Value = 35.5689531294117
$ws.Range("H9").Value = 16.1131547207898
$ws.Range("J9").Value = 10.59495259857779
$ws.Range("K9").Value = 19.58662681077378
$ws.Range("L9").Value = 9.818789910299119
$ws.Range("O9").Value = 25.41039108741667
$ws.Range("C10").Value = 14.03057897792292
$ws.Range("D10").Value = 9.677296113578411
$ws.Range("E10").Value = 14.58844324681213
$ws.Range("F10").Value = 33.61483566802742
$ws.Range("G10").Value = 35.22158807903214
$ws.Range("H10").Value = 15.99147403898007
$ws.Range("J10").Value = 10.56251930672673
$ws.Range("K10").Value = 20.64947570989334
$ws.Range("L10").Value = 9.806634786038588
$ws.Range("O10").Value = 25.1875529056669
$ws.Range("C11").Value = 14.08624771497803
$ws.Range("D11").Value = 9.724052905415057
$ws.Range("E11").Value = 14.60642063677027
$ws.Range("F11").Value = 33.52705708537945
$ws.Range("G11").Value = 35.07756006538823
$ws.Range("H11").Value = 15.93937748718533
$ws.Range("J11").Value = 10.54902752871484
$ws.Range("K11").Value = 21.11433878739299
$ws.Range("L11").Value = 9.801938455429624
$ws.Range("O11").Value = 25.09299718331559
$ws.Range("C12").Value = 14.10766898153631
$ws.Range("D12").Value = 9.741839351372107
$ws.Range("E12").Value = 14.61353384832825
$ws.Range("F12").Value = 33.49523668137394
$ws.Range("G12").Value = 35.02505049181858
$ws.Range("H12").Value = 15.92011814910214
$ws.Range("J12").Value = 10.54409971611568
$ws.Range("K12").Value = 21.28757822952999
$ws.Range("L12").Value = 9.800279371504447
$ws.Range("O12").Value = 25.05817452862486
$ws.Range("C13").Value = 14.10304059669027
$ws.Range("D13").Value = 9.738005299886856
$ws.Range("E13").Value = 14.6119883604644
$ws.Range("F13").Value = 33.50202652190541
$ws.Range("G13").Value = 35.03626874849722
$ws.Range("H13").Value = 15.92424515201128
$ws.Range("J13").Value = 10.54515295300251
$ws.Range("K13").Value = 21.25039395235967
$ws.Range("L13").Value = 9.800631386171538
$ws.Range("O13").Value = 25.06563040779348
$ws.Range("C14").Value = 14.08800330656129
$ws.Range("D14").Value = 9.72551464989318
$ws.Range("E14").Value = 14.60699974236631
$ws.Range("F14").Value = 33.52441070967648
$ws.Range("G14").Value = 35.07319926371051
$ws.Range("H14").Value = 15.93778361919443
$ws.Range("J14").Value = 10.54861848411096
$ws.Range("K14").Value = 21.12864790251089
$ws.Range("L14").Value = 9.80179957316013
$ws.Range("O14").Value = 25.09011256960905
$ws.Range("C15").Value = 14.07883649217961
$ws.Range("D15").Value = 9.717873962802351
$ws.Range("E15").Value = 14.60398375284811
$ws.Range("F15").Value = 33.5383067533371
$ws.Range("G15").Value = 35.09608529363825
$ws.Range("H15").Value = 15.9461373399185
$ws.Range("J15").Value = 10.55076481509265
$ws.Range("K15").Value = 21.05370790769501
$ws.Range("L15").Value = 9.802530645178395
$ws.Range("O15").Value = 25.10523679711816
$ws.Range("C16").Value = 14.02698932643922
$ws.Range("D16").Value = 9.674252299463385
$ws.Range("E16").Value = 14.58731134184746
$ws.Range("F16").Value = 33.62077041306199
$ws.Range("G16").Value = 35.23128378175068
$ws.Range("H16").Value = 15.99494421836143
$ws.Range("J16").Value = 10.56342640236626
$ws.Range("K16").Value = 20.61871071845322
$ws.Range("L16").Value = 9.806958425190173
$ws.Range("O16").Value = 25.19386977687884
$ws.Range("C17").Value = 13.99580434346254
$ws.Range("D17").Value = 9.64764656887227
$ws.Range("E17").Value = 14.57763151301866
$ws.Range("F17").Value = 33.67387909750307
$ws.Range("G17").Value = 35.31782031315844
$ws.Range("H17").Value = 16.02571987900189
$ws.Range("J17").Value = 10.57151697697231
$ws.Range("K17").Value = 20.34699522333731
$ws.Range("L17").Value = 9.809887759677899
$ws.Range("O17").Value = 25.24999113619062
$ws.Range("C18").Value = 13.97810032607401
$ws.Range("D18").Value = 9.632403835467072
$ws.Range("E18").Value = 14.57226651455051
$ws.Range("F18").Value = 33.70534934983366
$ws.Range("G18").Value = 35.36890950791029
$ws.Range("H18").Value = 16.04372773679832
$ws.Range("J18").Value = 10.57628929133137
$ws.Range("K18").Value = 20.18896373292787
$ws.Range("L18").Value = 9.811651079903749
$ws.Range("O18").Value = 25.28291184427304
$ws.Range("C19").Value = 13.97214643118018
$ws.Range("D19").Value = 9.627253556165579
$ws.Range("E19").Value = 14.57048491723701
$ws.Range("F19").Value = 33.71616301601284
$ws.Range("G19").Value = 35.38643282400336
$ws.Range("H19").Value = 16.04987752639483
$ws.Range("J19").Value = 10.57792553274012
$ws.Range("K19").Value = 20.1351605309915
$ws.Range("L19").Value = 9.812261597333961
$ws.Range("O19").Value = 25.29416824622325
$ws.Range("C20").Value = 13.99910004493352
$ws.Range("D20").Value = 9.650472640396034
$ws.Range("E20").Value = 14.57864100591613
$ws.Range("F20").Value = 33.66812994733611
$ws.Range("G20").Value = 35.30847204778257
$ws.Range("H20").Value = 16.02241203324656
$ws.Range("J20").Value = 10.57064342510308
$ws.Range("K20").Value = 20.37610158571812
$ws.Range("L20").Value = 9.809567811783262
$ws.Range("O20").Value = 25.24395053534404
$ws.Range("C21").Value = 14.0924109869393
$ws.Range("D21").Value = 9.729181350152226
$ws.Range("E21").Value = 14.60845675722906
$ws.Range("F21").Value = 33.51779734640385
$ws.Range("G21").Value = 35.06229661162821
$ws.Range("H21").Value = 15.93379432654982
$ws.Range("J21").Value = 10.547595657348
$ws.Range("K21").Value = 21.16448433979355
$ws.Range("L21").Value = 9.801453214064763
$ws.Range("O21").Value = 25.08289484254158
$ws.Range("C22").Value = 14.15537441034561
$ws.Range("D22").Value = 9.781087356858492
$ws.Range("E22").Value = 14.62972230689491
$ws.Range("F22").Value = 33.42782340937664
$ws.Range("G22").Value = 34.91325214233228
$ws.Range("H22").Value = 15.87860840802465
$ws.Range("J22").Value = 10.53358882692454
$ws.Range("K22").Value = 21.66340895341562
$ws.Range("L22").Value = 9.796845080610996
$ws.Range("O22").Value = 24.98337026411101
$ws.Range("C23").Value = 14.12159316135456
$ws.Range("D23").Value = 9.753345002773601
$ws.Range("E23").Value = 14.61821092750588
$ws.Range("F23").Value = 33.47508429400987
$ws.Range("G23").Value = 34.99170983080801
$ws.Range("H23").Value = 15.90781221796397
$ws.Range("J23").Value = 10.54096798678841
$ws.Range("K23").Value = 21.39865168394736
$ws.Range("L23").Value = 9.799241077310947
$ws.Range("O23").Value = 25.03596245421131
$ws.Range("C24").Value = 13.99760935754524
$ws.Range("D24").Value = 9.649194806808257
$ws.Range("E24").Value = 14.57818399073927
$ws.Range("F24").Value = 33.67072621885722
$ws.Range("G24").Value = 35.31269423155138
$ws.Range("H24").Value = 16.02390653058816
$ws.Range("J24").Value = 10.57103798117621
$ws.Range("K24").Value = 20.36294825506348
$ws.Range("L24").Value = 9.809712213470089
$ws.Range("O24").Value = 25.24667944846587
$ws.Range("C25").Value = 13.87282066211233
$ws.Range("D25").Value = 9.539192863611293
$ws.Range("E25").Value = 14.54276562189806
$ws.Range("F25").Value = 33.91703400245997
$ws.Range("G25").Value = 35.70915348273748
$ws.Range("H25").Value = 16.1608514441605
$ws.Range("J25").Value = 10.59495259857779
$ws.Range("K25").Value = 19.58662681077378
$ws.Range("L25").Value = 9.818789910299119
$ws.Range("O25").Value = 25.49847246677509
